$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (Spanish labels -> clean English field names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case municipality/state name corrections ---
# Connector words (de/del/el/y/los/las/la) capitalized to De/Del/El/Y/Los/Las/La,
# plus one-off fix: MonteMorelos -> Montemorelos
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B14").Value = "Playas De Rosarito"
$ws.Range("B34").Value = "Amatenango De La Frontera"
$ws.Range("B37").Value = "Bejucal De Ocampo"
$ws.Range("B39").Value = "Benemérito De Las Américas"
$ws.Range("B47").Value = "Chiapa De Corzo"
$ws.Range("B52").Value = "Comitán De Domínguez"
$ws.Range("B73").Value = "Marqués De Comillas"
$ws.Range("B78").Value = "Ocozocoautla De Espinosa"
$ws.Range("B86").Value = "Salto De Agua"
$ws.Range("B87").Value = "San Cristóbal De Las Casas"
$ws.Range("B120").Value = "Coyame Del Sotol"
$ws.Range("B128").Value = "Guadalupe Y Calvo"
$ws.Range("B130").Value = "Hidalgo Del Parral"
$ws.Range("B145").Value = "San Francisco De Borja"
$ws.Range("B146").Value = "San Francisco Del Oro"
$ws.Range("B171").Value = "San Juan De Sabinas"
$ws.Range("B187").Value = "Villa De Álvarez"
$ws.Range("A189").Value = "Ciudad De México"
$ws.Range("B193").Value = "Cuajimalpa De Morelos"
$ws.Range("B208").Value = "Coneto De Comonfort"
$ws.Range("B222").Value = "Nombre De Dios"
$ws.Range("B226").Value = "Pánuco De Coronado"
$ws.Range("B233").Value = "San Juan De Guadalupe"
$ws.Range("B234").Value = "San Juan Del Río"
$ws.Range("B235").Value = "San Pedro Del Gallo"
$ws.Range("A245").Value = "Estado De México"
$ws.Range("B245").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B248").Value = "Almoloya De Alquisiras"
$ws.Range("B249").Value = "Almoloya De Juárez"
$ws.Range("B250").Value = "Almoloya Del Río"
$ws.Range("B257").Value = "Atizapán De Zaragoza"
$ws.Range("B263").Value = "Chapa De Mota"
$ws.Range("B267").Value = "Coacalco De Berriozábal"
$ws.Range("B274").Value = "Ecatepec De Morelos"
$ws.Range("B282").Value = "Ixtapan De La Sal"
$ws.Range("B283").Value = "Ixtapan Del Oro"
$ws.Range("B298").Value = "Naucalpan De Juárez"
$ws.Range("B308").Value = "San Felipe Del Progreso"
$ws.Range("B309").Value = "San Martín De Las Pirámides"
$ws.Range("B311").Value = "San Simón De Guerrero"
$ws.Range("B313").Value = "Soyaniquilpan De Juárez"
$ws.Range("B323").Value = "Tenango Del Aire"
$ws.Range("B324").Value = "Tenango Del Valle"
$ws.Range("B337").Value = "Tlalnepantla De Baz"
$ws.Range("B343").Value = "Valle De Bravo"
$ws.Range("B344").Value = "Valle De Chalco Solidaridad"
$ws.Range("B345").Value = "Villa De Allende"
$ws.Range("B346").Value = "Villa Del Carbón"
$ws.Range("B360").Value = "Apaseo El Alto"
$ws.Range("B361").Value = "Apaseo El Grande"
$ws.Range("B369").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B373").Value = "Jaral Del Progreso"
$ws.Range("B381").Value = "Purísima Del Rincón"
$ws.Range("B385").Value = "San Diego De La Unión"
$ws.Range("B387").Value = "San Francisco Del Rincón"
$ws.Range("B389").Value = "San Luis De La Paz"
$ws.Range("B390").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B392").Value = "Silao De La Victoria"
$ws.Range("B397").Value = "Valle De Santiago"
$ws.Range("B403").Value = "Acapulco De Juárez"
$ws.Range("B406").Value = "Ajuchitlán Del Progreso"
$ws.Range("B407").Value = "Alcozauca De Guerrero"
$ws.Range("B411").Value = "Atenango Del Río"
$ws.Range("B412").Value = "Atlamajalcingo Del Monte"
$ws.Range("B414").Value = "Atoyac De Álvarez"
$ws.Range("B415").Value = "Ayutla De Los Libres"
$ws.Range("B418").Value = "Buenavista De Cuéllar"
$ws.Range("B419").Value = "Chilapa De Álvarez"
$ws.Range("B420").Value = "Chilpancingo De Los Bravo"
$ws.Range("B421").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B426").Value = "Coyuca De Benítez"
$ws.Range("B427").Value = "Coyuca De Catalán"
$ws.Range("B431").Value = "Cuetzala Del Progreso"
$ws.Range("B432").Value = "Cutzamala De Pinzón"
$ws.Range("B438").Value = "Huitzuco De Los Figueroa"
$ws.Range("B439").Value = "Iguala De La Independencia"
$ws.Range("B441").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B442").Value = "Zihuatanejo De Azueta"
$ws.Range("B444").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B447").Value = "Mártir De Cuilapan"
$ws.Range("B460").Value = "Taxco De Alarcón"
$ws.Range("B462").Value = "Técpan De Galeana"
$ws.Range("B464").Value = "Tepecoacuilco De Trujano"
$ws.Range("B466").Value = "Tixtla De Guerrero"
$ws.Range("B469").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B470").Value = "Tlapa De Comonfort"
$ws.Range("B482").Value = "Agua Blanca De Iturbide"
$ws.Range("B488").Value = "Atotonilco De Tula"
$ws.Range("B489").Value = "Atotonilco El Grande"
$ws.Range("B495").Value = "Cuautepec De Hinojosa"
$ws.Range("B501").Value = "Huasca De Ocampo"
$ws.Range("B505").Value = "Huejutla De Reyes"
$ws.Range("B508").Value = "Jacala De Ledezma"
$ws.Range("B513").Value = "Mineral Del Chico"
$ws.Range("B514").Value = "Mineral Del Monte"
$ws.Range("B515").Value = "Mixquiahuala De Juárez"
$ws.Range("B516").Value = "Molango De Escamilla"
$ws.Range("B518").Value = "Nopala De Villagrán"
$ws.Range("B519").Value = "Omitlán De Juárez"
$ws.Range("B520").Value = "Pachuca De Soto"
$ws.Range("B523").Value = "Progreso De Obregón"
$ws.Range("B528").Value = "Santiago De Anaya"
$ws.Range("B529").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B533").Value = "Tenango De Doria"
$ws.Range("B535").Value = "Tepehuacán De Guerrero"
$ws.Range("B536").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B539").Value = "Tezontepec De Aldama"
$ws.Range("B547").Value = "Tula De Allende"
$ws.Range("B548").Value = "Tulancingo De Bravo"
$ws.Range("B549").Value = "Villa De Tezontepec"
$ws.Range("B551").Value = "Zacualtipán De Ángeles"
$ws.Range("B556").Value = "Acatlán De Juárez"
$ws.Range("B557").Value = "Ahualulco De Mercado"
$ws.Range("B562").Value = "Atemajac De Brizuela"
$ws.Range("B565").Value = "Atotonilco El Alto"
$ws.Range("B567").Value = "Autlán De Navarro"
$ws.Range("B572").Value = "Cañadas De Obregón"
$ws.Range("B578").Value = "Concepción De Buenos Aires"
$ws.Range("B579").Value = "Cuautitlán De García Barragán"
$ws.Range("B588").Value = "Encarnación De Díaz"
$ws.Range("B595").Value = "Huejuquilla El Alto"
$ws.Range("B596").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B597").Value = "Ixtlahuacán Del Río"
$ws.Range("B601").Value = "Jilotlán De Los Dolores"
$ws.Range("B607").Value = "La Manzanilla De La Paz"
$ws.Range("B608").Value = "Lagos De Moreno"
$ws.Range("B615").Value = "Ojuelos De Jalisco"
$ws.Range("B620").Value = "San Cristóbal De La Barranca"
$ws.Range("B621").Value = "San Diego De Alejandría"
$ws.Range("B623").Value = "San Juan De Los Lagos"
$ws.Range("B626").Value = "San Martín De Bolaños"
$ws.Range("B628").Value = "San Miguel El Alto"
$ws.Range("B629").Value = "San Sebastián Del Oeste"
$ws.Range("B630").Value = "Santa María De Los Ángeles"
$ws.Range("B631").Value = "Santa María Del Oro"
$ws.Range("B634").Value = "Talpa De Allende"
$ws.Range("B635").Value = "Tamazula De Gordiano"
$ws.Range("B637").Value = "Techaluta De Montenegro"
$ws.Range("B641").Value = "Teocuitatlán De Corona"
$ws.Range("B642").Value = "Tepatitlán De Morelos"
$ws.Range("B645").Value = "Tizapán El Alto"
$ws.Range("B646").Value = "Tlajomulco De Zúñiga"
$ws.Range("B658").Value = "Unión De San Antonio"
$ws.Range("B659").Value = "Unión De Tula"
$ws.Range("B660").Value = "Valle De Guadalupe"
$ws.Range("B661").Value = "Valle De Juárez"
$ws.Range("B666").Value = "Yahualica De González Gallo"
$ws.Range("B667").Value = "Zacoalco De Torres"
$ws.Range("B670").Value = "Zapotitlán De Vadillo"
$ws.Range("B671").Value = "Zapotlán Del Rey"
$ws.Range("B672").Value = "Zapotlán El Grande"
$ws.Range("B698").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B700").Value = "Cojumatlán De Régules"
$ws.Range("B767").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B793").Value = "Coatlán Del Río"
$ws.Range("B801").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B805").Value = "Puente De Ixtla"
$ws.Range("B811").Value = "Tetela Del Volcán"
$ws.Range("B812").Value = "Tlaltizapán De Zapata"
$ws.Range("B820").Value = "Zacualpan De Amilpas"
$ws.Range("B824").Value = "Amatlán De Cañas"
$ws.Range("B825").Value = "Bahía De Banderas"
$ws.Range("B827").Value = "Ixtlán Del Río"
$ws.Range("B834").Value = "Santa María Del Oro"
$ws.Range("B858").Value = "Mier Y Noriega"
$ws.Range("B859").Value = "Montemorelos"
$ws.Range("B865").Value = "San Nicolás De Los Garza"
$ws.Range("B871").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B878").Value = "Capulálpam De Méndez"
$ws.Range("B880").Value = "Chalcatongo De Hidalgo"
$ws.Range("B881").Value = "Ciénega De Zimatlán"
$ws.Range("B884").Value = "Coicoyán De Las Flores"
$ws.Range("B887").Value = "Constancia Del Rosario"
$ws.Range("B890").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B892").Value = "Guevea De Humboldt"
$ws.Range("B893").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B894").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B895").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B896").Value = "Ixtlán De Juárez"
$ws.Range("B897").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B906").Value = "Mariscala De Juárez"
$ws.Range("B908").Value = "Mazatlán Villa De Flores"
$ws.Range("B910").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B911").Value = "Mixistlán De La Reforma"
$ws.Range("B914").Value = "Nejapa De Madero"
$ws.Range("B916").Value = "Oaxaca De Juárez"
$ws.Range("B917").Value = "Ocotlán De Morelos"
$ws.Range("B918").Value = "Pinotepa De Don Luis"
$ws.Range("B920").Value = "Putla Villa De Guerrero"
$ws.Range("B921").Value = "Rojas De Cuauhtémoc"
$ws.Range("B926").Value = "San Agustín De Las Juntas"
$ws.Range("B938").Value = "San Antonino El Alto"
$ws.Range("B943").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B953").Value = "San Dionisio Del Mar"
$ws.Range("B956").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B975").Value = "San José Del Peñasco"
$ws.Range("B981").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B1014").Value = "San Martín De Los Cansecos"
$ws.Range("B1019").Value = "San Mateo Del Mar"
$ws.Range("B1029").Value = "San Miguel Del Puerto"
$ws.Range("B1030").Value = "San Miguel El Grande"
$ws.Range("B1048").Value = "San Pablo Villa De Mitla"
$ws.Range("B1053").Value = "San Pedro El Alto"
$ws.Range("B1069").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1070").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B1071").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B1090").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B1094").Value = "Santa Inés De Zaragoza"
$ws.Range("B1095").Value = "Santa Inés Del Monte"
$ws.Range("B1097").Value = "Santa Lucía Del Camino"
$ws.Range("B1111").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1152").Value = "Santo Domingo De Morelos"
$ws.Range("B1164").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1165").Value = "Tataltepec De Valdés"
$ws.Range("B1166").Value = "Teotitlán De Flores Magón"
$ws.Range("B1167").Value = "Teotitlán Del Valle"
$ws.Range("B1169").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B1170").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1171").Value = "Tlacolula De Matamoros"
$ws.Range("B1172").Value = "Totontepec Villa De Morelos"
$ws.Range("B1174").Value = "Villa De Etla"
$ws.Range("B1175").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1176").Value = "Villa De Zaachila"
$ws.Range("B1178").Value = "Villa Sola De Vega"
$ws.Range("B1179").Value = "Villa Tejúpam De La Unión"
$ws.Range("B1182").Value = "Zapotitlán Del Río"
$ws.Range("B1185").Value = "Zimatlán De Álvarez"
$ws.Range("B1209").Value = "Ayotoxco De Guerrero"
$ws.Range("B1213").Value = "Chalchicomula De Sesma"
$ws.Range("B1222").Value = "Chila De La Sal"
$ws.Range("B1235").Value = "Cuayuca De Andrade"
$ws.Range("B1236").Value = "Cuetzalan Del Progreso"
$ws.Range("B1251").Value = "Huehuetlán El Chico"
$ws.Range("B1255").Value = "Huitzilan De Serdán"
$ws.Range("B1256").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1260").Value = "Izúcar De Matamoros"
$ws.Range("B1270").Value = "Los Reyes De Juárez"
$ws.Range("B1279").Value = "Palmar De Bravo"
$ws.Range("B1301").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1305").Value = "San Salvador El Seco"
$ws.Range("B1306").Value = "San Salvador El Verde"
$ws.Range("B1319").Value = "Tepanco De López"
$ws.Range("B1320").Value = "Tepango De Rodríguez"
$ws.Range("B1321").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1326").Value = "Tepexi De Rodríguez"
$ws.Range("B1328").Value = "Tetela De Ocampo"
$ws.Range("B1329").Value = "Teteles De Avila Castillo"
$ws.Range("B1334").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1344").Value = "Totoltepec De Guerrero"
$ws.Range("B1349").Value = "Xayacatlán De Bravo"
$ws.Range("B1354").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1369").Value = "Amealco De Bonfil"
$ws.Range("B1371").Value = "Cadereyta De Montes"
$ws.Range("B1376").Value = "Jalpan De Serra"
$ws.Range("B1377").Value = "Landa De Matamoros"
$ws.Range("B1380").Value = "Pinal De Amoles"
$ws.Range("B1382").Value = "San Juan Del Río"
$ws.Range("B1393").Value = "Armadillo De Los Infante"
$ws.Range("B1394").Value = "Axtla De Terrazas"
$ws.Range("B1399").Value = "Ciudad Del Maíz"
$ws.Range("B1409").Value = "Mexquitic De Carmona"
$ws.Range("B1415").Value = "San Ciro De Acosta"
$ws.Range("B1420").Value = "Santa María Del Río"
$ws.Range("B1422").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1432").Value = "Villa De Arista"
$ws.Range("B1433").Value = "Villa De Arriaga"
$ws.Range("B1434").Value = "Villa De Guadalupe"
$ws.Range("B1435").Value = "Villa De Ramos"
$ws.Range("B1436").Value = "Villa De Reyes"
$ws.Range("B1478").Value = "Nacozari De García"
$ws.Range("B1520").Value = "Soto La Marina"
$ws.Range("B1527").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1529").Value = "Amaxac De Guerrero"
$ws.Range("B1530").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1535").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1543").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1545").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1546").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1549").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1555").Value = "San Pablo Del Monte"
$ws.Range("B1561").Value = "Tepetitla De Lardizábal"
$ws.Range("B1564").Value = "Tetla De La Solidaridad"
$ws.Range("B1576").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1585").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1589").Value = "Amatlán De Los Reyes"
$ws.Range("B1601").Value = "Boca Del Río"
$ws.Range("B1603").Value = "Camarón De Tejeda"
$ws.Range("B1607").Value = "Castillo De Teayo"
$ws.Range("B1609").Value = "Cazones De Herrera"
$ws.Range("B1626").Value = "Cosamaloapan De Carpio"
$ws.Range("B1627").Value = "Cosautlán De Carvajal"
$ws.Range("B1643").Value = "Hueyapan De Ocampo"
$ws.Range("B1644").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1645").Value = "Ignacio De La Llave"
$ws.Range("B1649").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1650").Value = "Ixhuatlán De Madero"
$ws.Range("B1651").Value = "Ixhuatlán Del Café"
$ws.Range("B1652").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1662").Value = "Juchique De Ferrer"
$ws.Range("B1665").Value = "Landero Y Coss"
$ws.Range("B1668").Value = "Las Vigas De Ramírez"
$ws.Range("B1669").Value = "Lerdo De Tejada"
$ws.Range("B1673").Value = "Martínez De La Torre"
$ws.Range("B1675").Value = "Medellín De Bravo"
$ws.Range("B1679").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1690").Value = "Ozuluama De Mascareñas"
$ws.Range("B1694").Value = "Paso De Ovejas"
$ws.Range("B1695").Value = "Paso Del Macho"
$ws.Range("B1699").Value = "Poza Rica De Hidalgo"
$ws.Range("B1710").Value = "Sayula De Alemán"
$ws.Range("B1712").Value = "Soledad De Doblado"
$ws.Range("B1719").Value = "Tatahuicapan De Juárez"
$ws.Range("B1738").Value = "Tlacotepec De Mejía"
$ws.Range("B1751").Value = "Vega De Alatorre"
$ws.Range("B1761").Value = "Zozocolco De Hidalgo"
$ws.Range("B1775").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1777").Value = "Concepción Del Oro"
$ws.Range("B1779").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1790").Value = "Jiménez Del Teul"
$ws.Range("B1796").Value = "Mezquital Del Oro"
$ws.Range("B1801").Value = "Moyahua De Estrada"
$ws.Range("B1802").Value = "Nochistlán De Mejía"
$ws.Range("B1803").Value = "Noria De Ángeles"
$ws.Range("B1814").Value = "Teúl De González Ortega"
$ws.Range("B1815").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1817").Value = "Trinidad García De La Cadena"
$ws.Range("B1820").Value = "Villa De Cos"

# --- Percentage (column D) float re-normalization ---
# re-asserts the exact same ratios with the refreshed floating-point
# serialization produced by the cleaning re-run
$ws.Range("D36").Value = 0.0000960015360245764
$ws.Range("D51").Value = 0.0000960015360245764
$ws.Range("D74").Value = 0.0000960015360245764
$ws.Range("D77").Value = 0.0000960015360245764
$ws.Range("D88").Value = 0.0000960015360245764
$ws.Range("D97").Value = 0.0000960015360245764
$ws.Range("D109").Value = 0.0000960015360245764
$ws.Range("D115").Value = 0.0000960015360245764
$ws.Range("D149").Value = 0.0000960015360245764
$ws.Range("D178").Value = 0.0000960015360245764
$ws.Range("D192").Value = 0.009264148226371624
$ws.Range("D220").Value = 0.0000960015360245764
$ws.Range("D221").Value = 0.0000960015360245764
$ws.Range("D246").Value = 0.0000960015360245764
$ws.Range("D263").Value = 0.0000960015360245764
$ws.Range("D275").Value = 0.0000960015360245764
$ws.Range("D292").Value = 0.0000960015360245764
$ws.Range("D295").Value = 0.0000960015360245764
$ws.Range("D311").Value = 0.0000960015360245764
$ws.Range("D327").Value = 0.0000960015360245764
$ws.Range("D328").Value = 0.0000960015360245764
$ws.Range("D354").Value = 0.0000960015360245764
$ws.Range("D386").Value = 0.000944015104241668
$ws.Range("D402").Value = 0.09544152706443304
$ws.Range("D459").Value = 0.000944015104241668
$ws.Range("D497").Value = 0.0000960015360245764
$ws.Range("D530").Value = 0.0000960015360245764
$ws.Range("D546").Value = 0.0000960015360245764
$ws.Range("D566").Value = 0.0000960015360245764
$ws.Range("D597").Value = 0.0000960015360245764
$ws.Range("D602").Value = 0.0000960015360245764
$ws.Range("D603").Value = 0.0000960015360245764
$ws.Range("D694").Value = 0.0000960015360245764
$ws.Range("D714").Value = 0.0009120145922334756
$ws.Range("D789").Value = 0.0009120145922334756
$ws.Range("D832").Value = 0.0000960015360245764
$ws.Range("D854").Value = 0.0000960015360245764
$ws.Range("D905").Value = 0.0000960015360245764
$ws.Range("D914").Value = 0.0000960015360245764
$ws.Range("D980").Value = 0.0000960015360245764
$ws.Range("D997").Value = 0.0000960015360245764
$ws.Range("D1029").Value = 0.0000960015360245764
$ws.Range("D1050").Value = 0.0000960015360245764
$ws.Range("D1055").Value = 0.0000960015360245764
$ws.Range("D1066").Value = 0.0000960015360245764
$ws.Range("D1081").Value = 0.0000960015360245764
$ws.Range("D1121").Value = 0.0000960015360245764
$ws.Range("D1131").Value = 0.0000960015360245764
$ws.Range("D1165").Value = 0.0000960015360245764
$ws.Range("D1175").Value = 0.000960015360245764
$ws.Range("D1194").Value = 0.0000960015360245764
$ws.Range("D1211").Value = 0.0000960015360245764
$ws.Range("D1218").Value = 0.000944015104241668
$ws.Range("D1223").Value = 0.0000960015360245764
$ws.Range("D1229").Value = 0.0000960015360245764
$ws.Range("D1236").Value = 0.0000960015360245764
$ws.Range("D1267").Value = 0.0000960015360245764
$ws.Range("D1268").Value = 0.0000960015360245764
$ws.Range("D1277").Value = 0.0000960015360245764
$ws.Range("D1331").Value = 0.0000960015360245764
$ws.Range("D1335").Value = 0.000960015360245764
$ws.Range("D1375").Value = 0.0000960015360245764
$ws.Range("D1395").Value = 0.0000960015360245764
$ws.Range("D1423").Value = 0.0000960015360245764
$ws.Range("D1424").Value = 0.0000960015360245764
$ws.Range("D1429").Value = 0.0000960015360245764
$ws.Range("D1433").Value = 0.0000960015360245764
$ws.Range("D1454").Value = 0.0000960015360245764
$ws.Range("D1479").Value = 0.0000960015360245764
$ws.Range("D1528").Value = 0.0000960015360245764
$ws.Range("D1543").Value = 0.0000960015360245764
$ws.Range("D1555").Value = 0.0000960015360245764
$ws.Range("D1566").Value = 0.000960015360245764
$ws.Range("D1574").Value = 0.0000960015360245764
$ws.Range("D1590").Value = 0.0000960015360245764
$ws.Range("D1609").Value = 0.0000960015360245764
$ws.Range("D1615").Value = 0.0000960015360245764
$ws.Range("D1656").Value = 0.0000960015360245764
$ws.Range("D1674").Value = 0.0000960015360245764
$ws.Range("D1718").Value = 0.0000960015360245764
$ws.Range("D1792").Value = 0.0000960015360245764
$ws.Range("D1821").Value = 0.0000960015360245764

# --- Remove trailing metadata/footnote rows (1829:1833) ---
# sample-size / source / author / date footer lines are dropped;
# clearing collapses the used range back down to A1:D1827
$ws.Range("A1829:D1833").ClearContents()

Write-Output "edit applied"
